$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove rows whose term was dropped entirely (content merged into another row) ---
$ws.Rows.Item(49).Delete() | Out-Null

# --- Update remaining original rows 2-48 in place with revised term/other-terms/definition ---
$ws.Cells.Item(2,1).Value = "Study"
$ws.Cells.Item(2,2).ClearContents()
$ws.Cells.Item(2,3).Value = "A single funded research project resulting in one or more datasets to be used to answer a research question."
$ws.Cells.Item(3,1).Value = "Randomized controlled trial"
$ws.Cells.Item(3,2).Value = "RCT"
$ws.Cells.Item(3,3).Value = "A study design that randomly assigns participants to a control or treatment condition. In education research you often hear about two types of RCTs. The first being the Individual-Level Randomized Controlled Trial (I-RCT) in which individuals (such as students) are randomized directly to the treatment or control group. The second is a Cluster Randomized Controlled Trial (C-RCT), sometimes also called group-randomized, in which clusters of students (such as classrooms) are randomized."
$ws.Cells.Item(4,1).Value = "Subject"
$ws.Cells.Item(4,2).Value = "case, participant, site, record"
$ws.Cells.Item(4,3).Value = "A person or place participating in research and has one or more piece of data collected on them."
$ws.Cells.Item(5,1).Value = "Cohort"
$ws.Cells.Item(5,2).ClearContents()
$ws.Cells.Item(5,3).Value = "A group of participants recruited into the study at the same time."
$ws.Cells.Item(6,1).Value = "Wave"
$ws.Cells.Item(6,2).Value = "time period, time point, event, session"
$ws.Cells.Item(6,3).Value = "Intervals of data collection over time."
$ws.Cells.Item(7,1).Value = "Treatment"
$ws.Cells.Item(7,2).Value = "experiment"
$ws.Cells.Item(7,3).Value = "The individual or group receives the intervention."
$ws.Cells.Item(8,1).Value = "Control"
$ws.Cells.Item(8,2).Value = "business as usual"
$ws.Cells.Item(8,3).Value = "The individual or group does not receive the intervention."
$ws.Cells.Item(9,1).Value = "De-identified data"
$ws.Cells.Item(9,2).Value = "anonymized data, anonymization"
$ws.Cells.Item(9,3).Value = "Identifying information has been removed or distorted and the data can no longer be re-associated with the underlying individual (the linking key no longer exists). "
$ws.Cells.Item(10,1).Value = "Anonymous data"
$ws.Cells.Item(10,2).ClearContents()
$ws.Cells.Item(10,3).Value = "Identifying information was never collected. This data can not be linked across time or measures."
$ws.Cells.Item(11,1).Value = "Confidential data"
$ws.Cells.Item(11,2).Value = "pseudonymisation, coded data, indirectly identifiable"
$ws.Cells.Item(11,3).Value = "The status of this data is protected. Personally identifiable information (PII) in your data has been removed and names are replaced with a code and the only way to link the data back to an individual is through that code. The identifying code file (linking key) is stored separate from the research data."
$ws.Cells.Item(12,1).Value = "Identifiable data"
$ws.Cells.Item(12,2).ClearContents()
$ws.Cells.Item(12,3).Value = "Data that includes personally identifiable information."
$ws.Cells.Item(13,1).Value = "Personally identifiable information"
$ws.Cells.Item(13,2).Value = "PII"
$ws.Cells.Item(13,3).Value = "This includes direct identifiers (e.g., name and email), as well as indirect identifiers that, if combined with other variables, could identify a participant (e.g., full birthdate and county of residence). Under FERPA, additional PII, such as a district or school ID, should also be removed. Protected health identifiers (PHI) is a similar protected category of information. There are 18 HIPAA protected health identifiers that should be removed from data in order to meet the Safe Harbor de-identification method (e.g., name, email, address)."
$ws.Cells.Item(14,1).Value = "Unique participant identifier"
$ws.Cells.Item(14,2).Value = "study ID, site ID, unique identifier (UID), subject ID, participant code, record id"
$ws.Cells.Item(14,3).Value = "This is a numeric or alphanumeric identifier that is unique to every participant or site in order to create confidential and de-identified data. These identifiers allow researchers to link data across time or measure."
$ws.Cells.Item(15,1).Value = "Participant database"
$ws.Cells.Item(15,2).Value = "study roster, master list, master key, linking key, tracking database"
$ws.Cells.Item(15,3).Value = "This database, or spreadsheet, includes any identifiable information on your participants as well as their assigned study ID. It is your only own means of linking your confidential research study data to a participant’s true identity. It is also used to track data collected across time and measures as well as participant attrition."
$ws.Cells.Item(16,1).Value = "Attrition"
$ws.Cells.Item(16,2).ClearContents()
$ws.Cells.Item(16,3).Value = "The loss of study units from the sample, often seen in longitudinal studies."
$ws.Cells.Item(17,1).Value = "Longitudinal"
$ws.Cells.Item(17,2).ClearContents()
$ws.Cells.Item(17,3).Value = "Data is collected on participants over a period of time."
$ws.Cells.Item(18,1).Value = "Cross-sectional"
$ws.Cells.Item(18,2).ClearContents()
$ws.Cells.Item(18,3).Value = "Data is collected on participants for a single time point."
$ws.Cells.Item(19,1).Value = "File formats"
$ws.Cells.Item(19,2).ClearContents()
$ws.Cells.Item(19,3).Value = "Education research data is typically collected in one of three file formats: text( .txt, .pdf, .docx), tabular (.xlsx, .csv, .sav) , multimedia (.mpeg, .wav)."
$ws.Cells.Item(20,1).Value = "Dataset"
$ws.Cells.Item(20,2).Value = "dataframe, spreadsheet"
$ws.Cells.Item(20,3).Value = "A structured collection of data usually stored in tabular form. A research study usually produces one final dataset per entity/unit (e.g., teacher dataset, student dataset)."
$ws.Cells.Item(21,1).Value = "Raw data"
$ws.Cells.Item(21,2).Value = "primary, untouched"
$ws.Cells.Item(21,3).Value = "Unprocessed data collected directly from a source."
$ws.Cells.Item(22,1).Value = "Clean data"
$ws.Cells.Item(22,2).Value = "processed data"
$ws.Cells.Item(22,3).Value = "Raw data that has been manipulated or modified for the purposes of correcting and clarifying information."
$ws.Cells.Item(23,1).Value = "Database"
$ws.Cells.Item(23,2).Value = "relational database"
$ws.Cells.Item(23,3).Value = "An organized collection of related data stored in tables that can be linked together by a common identifier."
$ws.Cells.Item(24,1).Value = "Variable"
$ws.Cells.Item(24,2).Value = "column, field, question"
$ws.Cells.Item(24,3).Value = "Any phenomenon you are collecting information on/trying to measure. These variables will make up columns in your datasets or databases."
$ws.Cells.Item(25,1).Value = "Variable name"
$ws.Cells.Item(25,2).Value = "header"
$ws.Cells.Item(25,3).Value = "A shortened symbolic name given the variable in your data to represent the information it contains."
$ws.Cells.Item(26,1).Value = "Missing data"
$ws.Cells.Item(26,2).ClearContents()
$ws.Cells.Item(26,3).Value = "Occurs when there is no data stored in a variable for a particular observation/respondent."
$ws.Cells.Item(27,1).Value = "Directory"
$ws.Cells.Item(27,2).Value = "file structure, file tree"
$ws.Cells.Item(27,3).Value = "A cataloging structure for files and folders on your computer."
$ws.Cells.Item(28,1).Value = "Path"
$ws.Cells.Item(28,2).Value = "file path"
$ws.Cells.Item(28,3).Value = "A string of characters used to locate files in your directory system."
$ws.Cells.Item(29,1).Value = "Standardization"
$ws.Cells.Item(29,2).ClearContents()
$ws.Cells.Item(29,3).Value = "Developing a set of agreed upon technical standards and applying them within and across all research projects."
$ws.Cells.Item(30,1).Value = "Merge"
$ws.Cells.Item(30,2).Value = "join, link"
$ws.Cells.Item(30,3).Value = "Combining datasets together in a side by side manner (matching on one or more unique identifiers)."
$ws.Cells.Item(31,1).Value = "Append"
$ws.Cells.Item(31,2).ClearContents()
$ws.Cells.Item(31,3).Value = "Stacking datasets on top of each other (matching variables)."
$ws.Cells.Item(32,1).Value = "Syntax"
$ws.Cells.Item(32,2).Value = "code, program"
$ws.Cells.Item(32,3).Value = "Programming statements written in a text editor. The statements are machine-readable instructions processed by your computer."
$ws.Cells.Item(33,1).Value = "Qualitative data"
$ws.Cells.Item(33,2).ClearContents()
$ws.Cells.Item(33,3).Value = "Non-numeric data typically made up of text, images, video, or other artifacts."
$ws.Cells.Item(34,1).Value = "Quantitative data"
$ws.Cells.Item(34,2).ClearContents()
$ws.Cells.Item(34,3).Value = "Numerical data that can be analyzed with statistical methods."
$ws.Cells.Item(35,1).Value = "Reproducible"
$ws.Cells.Item(35,2).ClearContents()
$ws.Cells.Item(35,3).Value = "Being able to produce the same results using the same materials and procedures."
$ws.Cells.Item(36,1).Value = "Replicable"
$ws.Cells.Item(36,2).ClearContents()
$ws.Cells.Item(36,3).Value = "Being able to produce the same results if the same procedures are used with different materials."
$ws.Cells.Item(37,1).Value = "Observational data"
$ws.Cells.Item(37,2).ClearContents()
$ws.Cells.Item(37,3).Value = "Data collected from a study where researchers are observing the effect of an intervention without manipulating who is exposed to the intervention. This includes many formats that education researchers collect data with (e.g., survey, observation, assessment)."
$ws.Cells.Item(38,1).Value = "Experimental data"
$ws.Cells.Item(38,2).ClearContents()
$ws.Cells.Item(38,3).Value = "Data collected from a study where researchers randomly introduce an intervention and study the effects."
$ws.Cells.Item(39,1).Value = "Archive"
$ws.Cells.Item(39,2).ClearContents()
$ws.Cells.Item(39,3).Value = "The transfer of data to a facility, such as a repository, that preserves and stores data long-term."
$ws.Cells.Item(40,1).Value = "Privacy"
$ws.Cells.Item(40,2).ClearContents()
$ws.Cells.Item(40,3).Value = "Privacy concerns people, ensuring they are given control to the access of themselves and their information."
$ws.Cells.Item(41,1).Value = "Data"
$ws.Cells.Item(41,2).Value = "research data"
$ws.Cells.Item(41,3).Value = "The recorded factual material commonly accepted in the scientific community as necessary to validate research findings. (OMB Circular A-110)"
$ws.Cells.Item(42,1).Value = "Confidentiality"
$ws.Cells.Item(42,2).ClearContents()
$ws.Cells.Item(42,3).Value = "Confidentiality concerns data, ensuring participants agree to how their private and identifable information will be managed and disseminated."
$ws.Cells.Item(43,1).Value = "Simulation data"
$ws.Cells.Item(43,2).ClearContents()
$ws.Cells.Item(43,3).Value = "Data generated through imitations of a real-world process using computer models."
$ws.Cells.Item(44,1).Value = "Derived data"
$ws.Cells.Item(44,2).ClearContents()
$ws.Cells.Item(44,3).Value = "Data created through transformations of existing data (e.g., mean scores)."
$ws.Cells.Item(45,1).Value = "Extant data"
$ws.Cells.Item(45,2).Value = "secondary data"
$ws.Cells.Item(45,3).Value = "Existing data generated/collected by external organizations such as governments at an earlier point in time (e.g., administrative data)."
$ws.Cells.Item(46,1).Value = "Data type"
$ws.Cells.Item(46,2).Value = "measurement unit, variable format, variable class"
$ws.Cells.Item(46,3).Value = "A classification that specifies what types of values are contained in a variable and what kinds of operations can be performed on that variable. Examples of types include numeric, character, logical, or datetime."
$ws.Cells.Item(47,1).Value = "Research"
$ws.Cells.Item(47,2).ClearContents()
$ws.Cells.Item(47,3).Value = "The Common Rule (45 CFR 46) definition of research is a systematic investigation, including research development, testing, and evaluation, designed to develop or contribute to generalizable knowledge."
$ws.Cells.Item(48,1).Value = "Primary data"
$ws.Cells.Item(48,2).Value = "original data"
$ws.Cells.Item(48,3).Value = "First hand data that is generated/collected by the research team as part of the research study."

# --- Append 10 brand-new glossary terms starting at row 49 ---
$ws.Cells.Item(49,1).Value = "Direct identifiers"
$ws.Cells.Item(49,2).ClearContents()
$ws.Cells.Item(49,3).Value = "These identifiers can directly identify a participant and should always be removed from research study data. There should be no need to keep these identifiers for analysis (i.e. name, email, address)."
$ws.Cells.Item(50,1).Value = "FERPA"
$ws.Cells.Item(50,2).ClearContents()
$ws.Cells.Item(50,3).Value = "The Family Educational Rights and Privacy Act is a federal law governing the disclosure of personally identifiable information in education records (e.g., name, address, DOB). The law applies to all public elementary and secondary schools, as well as post-secondary institutions."
$ws.Cells.Item(51,1).Value = "HIPAA"
$ws.Cells.Item(51,2).ClearContents()
$ws.Cells.Item(51,3).Value = "The Health Insurance Portability and Accountability Act is a federal law covering the protection of sensitive health information."
$ws.Cells.Item(52,1).Value = "Indirect identifiers"
$ws.Cells.Item(52,2).ClearContents()
$ws.Cells.Item(52,3).Value = "Even though these identifiers are not necessarily uniquely tied to one individual (i.e., birthdate or place of birth), if combined, this information could indirectly identify a participant. Therefore this information should be managed before publicly sharing data."
$ws.Cells.Item(53,1).Value = "Instrument"
$ws.Cells.Item(53,2).ClearContents()
$ws.Cells.Item(53,3).Value = "A mechanism designed to collect original data (e.g., observation form, questionnaire, assessment)"
$ws.Cells.Item(54,1).Value = "Measure"
$ws.Cells.Item(54,2).ClearContents()
$ws.Cells.Item(54,3).Value = "In this book, I use the term measure broadly to refer to a collection of items used to measure an outcome (e.g., an existing scale, an existing academic assessment)."
$ws.Cells.Item(55,1).Value = "Private data"
$ws.Cells.Item(55,2).ClearContents()
$ws.Cells.Item(55,3).Value = "Highly restricted data with limited access (i.e., passwords)."
$ws.Cells.Item(56,1).Value = "Scale"
$ws.Cells.Item(56,2).ClearContents()
$ws.Cells.Item(56,3).Value = "Similar to the term `"measure`", this is a collection of items used to measure an outcome. However, I typically use this term to more specifically refer to questionnaires that have had psychometric properties assessed. Scales may also be made up of subscales (i.e., groupings of items)."
$ws.Cells.Item(57,1).Value = "Sensitive data"
$ws.Cells.Item(57,2).ClearContents()
$ws.Cells.Item(57,3).Value = "Private information that could cause harm and should be protected from unwarranted disclosure."
$ws.Cells.Item(58,1).Value = "Tool"
$ws.Cells.Item(58,2).ClearContents()
$ws.Cells.Item(58,3).Value = "A means used to collect data using an instrument (e.g., a paper form, an online survey platform)"

# --- Sort the whole glossary range alphabetically by Term (column A) ---
$sortRange = $ws.Range("A2:C58")
$keyRange = $ws.Range("A2:A58")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Match the saved view/selection state ---
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("C41").Select() | Out-Null
